$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add the new header title in C1
$ws.Range("C1").Value = "Name of Rmd file"

# Update the active selection to C1 (as in the saved sheet view)
$ws.Activate()
$ws.Range("C1").Select()
